$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "calendarios-regalos" rows (old rows 52-56). This shifts
#    the "invitaciones-papeleria" rows (old 57-64) up to become rows 52-59.
# ---------------------------------------------------------------------------
$ws.Range("A52:G56").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Append the new "calendarios-regalos" subcategory rows (16 rows) at the
#    bottom of the sheet, starting at row 60.
# ---------------------------------------------------------------------------
$newRows = @(
    @("albumes-fotos-familia", "Albumes Fotos Familia", "calendarios-regalos", "Subcategoría de Albumes Fotos Familia", "/media/subcategory_images/calendarios_regalos/albumes-fotos-familia.jpg", 1, "card"),
    @("almohadas-mascota-nuevo", "Almohadas Mascota", "calendarios-regalos", "Subcategoría de Almohadas Mascota", "/media/subcategory_images/calendarios_regalos/almohadas-mascota-nuevo.jpg", 2, "card"),
    @("arte-pared-familia", "Arte Pared Familia", "calendarios-regalos", "Subcategoría de Arte Pared Familia", "/media/subcategory_images/calendarios_regalos/arte-pared-familia.jpg", 3, "card"),
    @("tazas-termos", "Bebidas Termos Familia", "calendarios-regalos", "Subcategoría de Bebidas Termos Familia", "/media/subcategory_images/calendarios_regalos/tazas-termos.jpg", 4, "card"),
    @("bolsas-tote-personalizadas", "Bolsas Tote Personalizadas", "calendarios-regalos", "Subcategoría de Bolsas Tote Personalizadas", "/media/subcategory_images/calendarios_regalos/bolsas-tote-personalizadas.jpg", 5, "card"),
    @("calendarios-familiares", "Calendarios Familia", "calendarios-regalos", "Subcategoría de Calendarios ", "/media/subcategory_images/calendarios_regalos/calendarios-familia.jpg", 6, "card"),
    @("cartas-juego-personalizadas-nuevo", "Cartas Juego Personalizadas", "calendarios-regalos", "Subcategoría de Cartas Juego Personalizadas", "/media/subcategory_images/calendarios_regalos/cartas-juego-personalizadas-nuevo.jpg", 7, "card"),
    @("cobijas-mantas-personalizadas", "Cobijas Mantas Personalizadas", "calendarios-regalos", "Subcategoría de Cobijas Mantas Personalizadas", "/media/subcategory_images/calendarios_regalos/cobijas-mantas-personalizadas.jpg", 8, "card"),
    @("decoracion-hogar-familia", "Decoracion Hogar Familia", "calendarios-regalos", "Subcategoría de Decoracion Hogar Familia", "/media/subcategory_images/calendarios_regalos/decoracion-hogar-familia.jpg", 9, "card"),
    @("impresiones-canvas-nuevo", "Impresiones Canvas", "calendarios-regalos", "Subcategoría de Impresiones Canvas", "/media/subcategory_images/calendarios_regalos/impresiones-canvas-nuevo.jpg", 10, "card"),
    @("mousepad-personalizado", "Mousepad Personalizado", "calendarios-regalos", "Subcategoría de Mousepad Personalizado", "/media/subcategory_images/calendarios_regalos/mousepad-personalizado.jpg", 11, "card"),
    @("papeleria-personal-familia", "Papeleria Personal Familia", "calendarios-regalos", "Subcategoría de Papeleria Personal Familia", "/media/subcategory_images/calendarios_regalos/papeleria-personal-familia.jpg", 12, "card"),
    @("plantilla-coordenadas-clasicas", "Plantilla Coordenadas Clasicas", "calendarios-regalos", "Subcategoría de Plantilla Coordenadas Clasicas", "/media/subcategory_images/calendarios_regalos/plantilla-coordenadas-clasicas.jpg", 13, "card"),
    @("plantilla-monograma-flores", "Plantilla Monograma Flores", "calendarios-regalos", "Subcategoría de Plantilla Monograma Flores", "/media/subcategory_images/calendarios_regalos/plantilla-monograma-flores.jpg", 14, "card"),
    @("servilletas-foil-nuevo", "Servilletas Foil", "calendarios-regalos", "Subcategoría de Servilletas Foil", "/media/subcategory_images/calendarios_regalos/servilletas-foil-nuevo.jpg", 15, "card"),
    @("tarjetas-notas-personalizadas", "Tarjetas Notas Personalizadas", "calendarios-regalos", "Subcategoría de Tarjetas Notas Personalizadas", "/media/subcategory_images/calendarios_regalos/tarjetas-notas-personalizadas.jpg", 16, "card"),
)

$startRow = 60
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 3. Update the hidden _FilterDatabase defined name so it reflects the new
#    (smaller) table bounds used for the old filter selection.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "subcategories_complete!_FilterDatabase") {
        $n.RefersTo = "=subcategories_complete!`$A`$1:`$G`$52"
    }
}

# ---------------------------------------------------------------------------
# 4. Restore the selection/active cell & scroll position to match the
#    author's final view of the sheet (row 66, template selection).
# ---------------------------------------------------------------------------
$ws.Range("A66").Select()

Write-Host "edit applied"
